$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 1553.5
$ws.Range("I38").Value = 285.6
$ws.Range("K38").Value = 856.8000000000001
$ws.Range("M38").Value = -484.8000000000001
# Row 41
$ws.Range("H41").Value = 916.1667
$ws.Range("I41").Value = 919.8
$ws.Range("K41").Value = 919.8
$ws.Range("M41").Value = -479.8
# Row 70
$ws.Range("H70").Value = 5294.7036
$ws.Range("I70").Value = 4712.6
$ws.Range("J70").Value = 5637.1177
$ws.Range("K70").Value = 14137.8
$ws.Range("L70").Value = 16911.3531
$ws.Range("M70").Value = -13867.8
$ws.Range("N70").Value = -17451.3531
# Row 73
$ws.Range("H73").Value = 5294.7036
$ws.Range("I73").Value = 4712.6
$ws.Range("J73").Value = 5637.1177
$ws.Range("K73").Value = 14137.8
$ws.Range("L73").Value = 16911.3531
$ws.Range("M73").Value = -13201.8
$ws.Range("N73").Value = -18783.3531
# Row 111
$ws.Range("H111").Value = 4478.375
$ws.Range("I111").Value = 4825
$ws.Range("J111").Value = 4131.75
$ws.Range("K111").Value = 14475
$ws.Range("L111").Value = 12395.25
$ws.Range("M111").Value = -11408
$ws.Range("N111").Value = -18529.25
# Row 113
$ws.Range("H113").Value = 5271.636
$ws.Range("I113").Value = 3248.625
$ws.Range("J113").Value = 10666.333
$ws.Range("K113").Value = 3248.625
$ws.Range("L113").Value = 10666.333
$ws.Range("M113").Value = 5.375
$ws.Range("N113").Value = -17174.333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2508.3635
$ws.Range("I45").Value = 2335.375
$ws.Range("J45").Value = 2969.6667
$ws.Range("K45").Value = 2335.375
$ws.Range("L45").Value = 2969.6667
$ws.Range("M45").Value = -1958.375
$ws.Range("N45").Value = -3723.6667
# Row 63
$ws.Range("H63").Value = 13078.6875
$ws.Range("I63").Value = 11489.154
$ws.Range("K63").Value = 11489.154
$ws.Range("M63").Value = -10803.154
# Row 66
$ws.Range("H66").Value = 13078.6875
$ws.Range("I66").Value = 11489.154
$ws.Range("K66").Value = 57445.77
$ws.Range("M66").Value = -54013.77
# Row 132
$ws.Range("H132").Value = 3188
$ws.Range("I132").Value = 3086.1072
$ws.Range("J132").Value = 3901.25
$ws.Range("K132").Value = 9258.321599999999
$ws.Range("L132").Value = 11703.75
$ws.Range("M132").Value = -6728.321599999999
$ws.Range("N132").Value = -16763.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 6
$ws.Range("H6").Value = 9045
$ws.Range("J6").Value = 9045
$ws.Range("L6").Value = 9045
$ws.Range("N6").Value = -9271
# Row 62
$ws.Range("H62").Value = 60000
$ws.Range("J62").Value = 60000
$ws.Range("L62").Value = 60000
$ws.Range("N62").Value = -61372
# Row 65
$ws.Range("H65").Value = 60000
$ws.Range("J65").Value = 60000
$ws.Range("L65").Value = 180000
$ws.Range("N65").Value = -186864
# Row 86
$ws.Range("H86").Value = 4580.421
$ws.Range("I86").Value = 2618.2666
$ws.Range("K86").Value = 2618.2666
$ws.Range("M86").Value = -1495.2666
# Row 89
$ws.Range("H89").Value = 4580.421
$ws.Range("I89").Value = 2618.2666
$ws.Range("K89").Value = 13091.333
$ws.Range("M89").Value = -7475.332999999999
# Row 105
$ws.Range("H105").Value = 5775.269
$ws.Range("I105").Value = 5774.143
$ws.Range("K105").Value = 5774.143
$ws.Range("M105").Value = -4027.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 734.8570999999999
$ws.Range("I16").Value = 734.8570999999999
$ws.Range("K16").Value = 734.8570999999999
$ws.Range("M16").Value = -447.8570999999999
# Row 99
$ws.Range("H99").Value = 6169.8
$ws.Range("I99").Value = 3962.25
$ws.Range("K99").Value = 3962.25
$ws.Range("M99").Value = -2464.25
# Row 113
$ws.Range("H113").Value = 734.8570999999999
$ws.Range("I113").Value = 734.8570999999999
$ws.Range("K113").Value = 734.8570999999999
$ws.Range("M113").Value = 1435.1429
# Row 122
$ws.Range("H122").Value = 3353.9375
$ws.Range("I122").Value = 3365.8462
$ws.Range("K122").Value = 10097.5386
$ws.Range("M122").Value = -7647.5386
# Row 126
$ws.Range("H126").Value = 6169.8
$ws.Range("I126").Value = 3962.25
$ws.Range("K126").Value = 11886.75
$ws.Range("M126").Value = -9416.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 796
$ws.Range("I14").Value = 796
$ws.Range("K14").Value = 2388
$ws.Range("M14").Value = -2215
# Row 23
$ws.Range("H23").Value = 323.4375
$ws.Range("I23").Value = 221.8
$ws.Range("J23").Value = 369.63635
$ws.Range("K23").Value = 665.4000000000001
$ws.Range("L23").Value = 1108.90905
$ws.Range("M23").Value = -430.4000000000001
$ws.Range("N23").Value = -1578.90905
# Row 33
$ws.Range("H33").Value = 170
$ws.Range("I33").Value = 169
$ws.Range("J33").Value = 172
$ws.Range("K33").Value = 1014
$ws.Range("L33").Value = 1032
$ws.Range("M33").Value = -731
$ws.Range("N33").Value = -1598
# Row 46
$ws.Range("H46").Value = 605.2
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 51
$ws.Range("H51").Value = 998
$ws.Range("I51").Value = 998
$ws.Range("K51").Value = 2994
$ws.Range("M51").Value = -2534
# Row 132
$ws.Range("H132").Value = 1954.7778
$ws.Range("I132").Value = 1927.5714
$ws.Range("J132").Value = 2050
$ws.Range("K132").Value = 17348.1426
$ws.Range("L132").Value = 18450
$ws.Range("M132").Value = -14818.1426
$ws.Range("N132").Value = -23510

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2157.842
$ws.Range("I7").Value = 2157.842
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2157.842
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2045.842
$ws.Range("N7").ClearContents()
# Row 16
$ws.Range("H16").Value = 350.66666
$ws.Range("I16").Value = 350.66666
$ws.Range("K16").Value = 350.66666
$ws.Range("M16").Value = -180.66666
# Row 118
$ws.Range("H118").Value = 43200
$ws.Range("J118").Value = 43200
$ws.Range("L118").Value = 43200
$ws.Range("N118").Value = -46514
# Row 122
$ws.Range("H122").Value = 2950.9167
$ws.Range("I122").Value = 2950.9167
$ws.Range("K122").Value = 8852.750100000001
$ws.Range("M122").Value = -6402.750100000001
# Row 126
$ws.Range("H126").Value = 2157.842
$ws.Range("I126").Value = 2157.842
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6473.526
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4003.526
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 103
$ws.Range("H103").Value = 30499.75
$ws.Range("J103").Value = 30499.75
$ws.Range("L103").Value = 30499.75
$ws.Range("N103").Value = -32843.75
